# Helper: force a value to be written as TEXT, even when it looks numeric
# (e.g. "000066", "0.70", "4.44"), then strip the temporary "@" number
# format back off so the cell ends up with no explicit style applied
# (matching the source data's plain, unstyled text cells).
function Set-TextCell($rng, $val) {
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. "总计" sheet: insert a new row for 2022-Q4 right after the header,
#    pushing 2021-Q2 / 2021-Q1 / 2020-Q4 down one row and renumbering
#    the index column (A).
# ------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

$summary.Rows.Item(2).Insert()
$summary.Range("B2:D2").ClearFormats()

# Give the new index cell (A2) the same style as the other index cells
# (bold/centered/bordered) by copying the format from A3.
$summary.Range("A3").Copy()
$summary.Range("A2").PasteSpecial(-4122)

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q4"
$summary.Range("C2").Value = 2
$summary.Range("D2").Value = 0.03

$summary.Range("A3").Value = 1
$summary.Range("A4").Value = 2
$summary.Range("A5").Value = 3

# ------------------------------------------------------------------
# 2. Insert a brand-new "2022-Q4" worksheet right before "2021-Q2" by
#    duplicating the "2021-Q2" sheet (so the header/index-column
#    styling comes across exactly), then replace its data with the
#    2022-Q4 fund-holding figures.
# ------------------------------------------------------------------
$template = $wb.Worksheets.Item("2021-Q2")
$template.Copy($template)
$q4 = $wb.Worksheets.Item("2021-Q2 (2)")
$q4.Name = "2022-Q4"

# The template sheet has 6 rows (header + 5 funds); 2022-Q4 only needs
# header + 2 funds, so drop the extra rows.
$q4.Rows("4:6").Delete()

$q4.Range("D1").Value = "基金规模"

$q4.Range("A2").Value = 0
Set-TextCell $q4.Range("B2") "000066"
Set-TextCell $q4.Range("C2") "诺安鸿鑫混合A"
Set-TextCell $q4.Range("D2") "0.70"
Set-TextCell $q4.Range("E2") "90.20"
Set-TextCell $q4.Range("F2") "4.44"
Set-TextCell $q4.Range("G2") "0.0311"
$q4.Range("H2").Value = 7

$q4.Range("A3").Value = 1
Set-TextCell $q4.Range("B3") "014498"
Set-TextCell $q4.Range("C3") "诺安鸿鑫混合C"
Set-TextCell $q4.Range("D3") "0.01"
Set-TextCell $q4.Range("E3") "90.20"
Set-TextCell $q4.Range("F3") "4.44"
Set-TextCell $q4.Range("G3") "0.0004"
$q4.Range("H3").Value = 7
